$p = $ppt.ActivePresentation

function Set-DatePlaceholderText($container, [string]$newText) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

# --- Update cached "datetimeFigureOut" date text 6/28/16 -> 6/30/16 ---
# Slide master
Set-DatePlaceholderText $p.SlideMaster "6/30/16"

# All slide layouts that carry their own date placeholder text
$layoutIndices = @(1,2,3,5,6,9,10,11)
foreach ($li in $layoutIndices) {
    $lay = $p.SlideMaster.CustomLayouts.Item($li)
    Set-DatePlaceholderText $lay "6/30/16"
}

# Handout master
Set-DatePlaceholderText $p.HandoutMaster "6/30/16"

# Notes master
Set-DatePlaceholderText $p.NotesMaster "6/30/16"

# --- Slide 22 table: bold the left-hand column text ---
$s = $p.Slides.Item(22)
$tbl = $s.Shapes.Item(2).Table
for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    $tbl.Cell($r, 1).Shape.TextFrame.TextRange.Font.Bold = $true
}
